$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 45132.4698562732
$ws.Range("E3").Value = 45132.4698563079
$ws.Range("E4").Value = 45132.4698563079
$ws.Range("E5").Value = 45132.4698563079
$ws.Range("E6").Value = 45132.4698563079
$ws.Range("E7").Value = 45132.4698563079
$ws.Range("E8").Value = 45132.4698563079
